$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The team renumbered the manual test cases: a new row (A19) was filled in
# with serial number 18, so every following row's "SL. No" shifts up by one.
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22

# Scroll the sheet so row 17 is at the top of the view, then select A24 as
# the active cell to match where editing left off.
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A24").Select()
